# Applies the LOQ4251.docx edit: convert the long run-on "Programa" (PT/EN)
# and "Bibliografia" paragraphs into one run per paragraph containing a
# <w:t>/<w:br/> pair for each logical line (matching the source document's
# existing style, e.g. the "Creditos-aula" paragraph).

$d = $word.ActiveDocument

function Get-ParagraphRangeByText($doc, [string]$needle) {
    $hit = $doc.Content
    $ok = $hit.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find paragraph containing: $needle"
    }
    return $hit.Paragraphs(1).Range
}

# 1) "Programa" section - Portuguese program paragraph
$ptRange = Get-ParagraphRangeByText $d "Programa em português1.Conceitos básicos"
$ptXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Programa em português</w:t><w:br/><w:t>1.Conceitos básicos de Química (2 horas)</w:t><w:br/><w:t>a.Estrutura Atômica</w:t><w:br/><w:t>b.Tabela Periódica</w:t><w:br/><w:t>c.Ligações Químicas</w:t><w:br/><w:t>2.Os estados físicos da matéria e suas propriedades peculiares (6 horas)</w:t><w:br/><w:t>a.O estado gasoso – pressão, relações PVT, gases ideais e reais</w:t><w:br/><w:t>b.O estado líquido – soluções, forças intermoleculares, viscosidade, tensão superficial, pressão de vapor, mudanças de fase</w:t><w:br/><w:t xml:space="preserve">c.O estado sólido – classificação dos sólidos (moleculares, reticulares, metálicos e iônicos) </w:t><w:br/><w:t>3.Reações químicas (8 horas)</w:t><w:br/><w:t>a.Tipos de reações (dupla-troca, oxirredução)</w:t><w:br/><w:t>b.Estequiometria em reações químicas (reagentes limitantes, pureza e rendimento)</w:t><w:br/><w:t>c.Energia e reações químicas</w:t><w:br/><w:t>d.Equilíbrio químico – soluções tampão</w:t><w:br/><w:t>e.Fundamentos de corrosão</w:t><w:br/><w:t>4.Noções de química orgânica (6 horas)</w:t><w:br/><w:t>a.Hidrocarbonetos e suas principais propriedades</w:t><w:br/><w:t>b.Combustíveis e combustão</w:t><w:br/><w:t>c.Polímeros</w:t><w:br/><w:t>5.Tecnologia Química aplicada (8 horas)</w:t><w:br/><w:t>a.Papel e celulose</w:t><w:br/><w:t>b.Açúcar e álcool</w:t><w:br/><w:t>c.Sabões e detergentes</w:t><w:br/><w:t>d.Petróleo e gás</w:t><w:br/><w:t xml:space="preserve">e.Gases industriais </w:t><w:br/><w:t>f.Produção de vidros e cimento</w:t></w:r></w:p>'
$ptRange.InsertXML($ptXml)

# 2) "Programa" section - English (italic) program paragraph
$enRange = Get-ParagraphRangeByText $d "1.Basic Concepts of Chemistrya.Atomic Structure"
$enXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>1.Basic Concepts of Chemistry</w:t><w:br/><w:t>a.Atomic Structure</w:t><w:br/><w:t>b.Periodic table</w:t><w:br/><w:t>c.Chemical bonds</w:t><w:br/><w:t>2.The physical states of matter and their peculiar properties</w:t><w:br/><w:t>a.The gaseous state – pressure, PVT relations, ideal and real gases</w:t><w:br/><w:t>b.The liquid state - solutions, intermolecular forces, viscosity, surface tension, vapor pressure, phase changes</w:t><w:br/><w:t>c.The solid state - classification of solids (molecular, reticular, metallic and ionic)</w:t><w:br/><w:t>3.Chemical reactions</w:t><w:br/><w:t>a.Types of reactions (double-exchange, oxy-reduction)</w:t><w:br/><w:t xml:space="preserve">b.Stoichiometry in chemical reactions (limiting reagents, purity and yield) </w:t><w:br/><w:t>c.Energy and chemical reactions</w:t><w:br/><w:t>d.Corrosion Fundamentals</w:t><w:br/><w:t>4.Notions of organic chemistry</w:t><w:br/><w:t>a.Hydrocarbons and their main properties</w:t><w:br/><w:t>b.Fuel and combustion</w:t><w:br/><w:t>c.Polymers</w:t><w:br/><w:t>5.Applied Chemistry Technology</w:t><w:br/><w:t>a.Paper and Cellulose</w:t><w:br/><w:t>b.Sugar and alcohol</w:t><w:br/><w:t>c.Soaps and detergents</w:t><w:br/><w:t>d.Oil and gas</w:t><w:br/><w:t>e.Industrial gases</w:t><w:br/><w:t>f.Glass and cement production</w:t></w:r></w:p>'
$enRange.InsertXML($enXml)

# 3) "Bibliografia" section paragraph
$bibRange = Get-ParagraphRangeByText $d "BROWN, T.L. et al. Química a ciência central"
$bibXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>BROWN, T.L. et al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007</w:t><w:br/><w:t>ATKINS, P. Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006</w:t><w:br/><w:t>KOTZ, J. C. et al. Química geral e reações químicas, 9ª Edição, São Paulo, Cengage Learning, 2015.</w:t><w:br/><w:t>TOLENTINO, N. M. C. Processos Químicos Industriais, 1ª Edição, São Paulo, Érica, 2015.</w:t></w:r></w:p>'
$bibRange.InsertXML($bibXml)

Write-Host "Applied Programa (PT/EN) and Bibliografia line-break reformatting."
